# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" / clrScheme "Office"      (currently wired to the notes master only)
#   ppt/theme/theme2.xml -> "Integral"     / clrScheme "Red Violet"  (wired to the slide master + presentation)
#
# The authored change swaps the *content* of the two theme parts (the "Office"
# colors move into theme1.xml, the "Red Violet"/Integral colors move into
# theme2.xml) while every relationship stays pointed at the same part names.
# Net visible effect on the slide master / presentation theme (the only theme
# this PowerPoint host's object model exposes for editing) is that its color
# scheme switches from the "Red Violet" Integral palette to the stock
# "Office" palette. We reproduce that by rewriting the 12 theme colors on the
# presentation's ThemeColorScheme to the Office defaults, in the fixed
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink.

function Hex-ToRgbLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p   = $ppt.ActivePresentation
$sm  = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $tcs.Item($i).RGB = Hex-ToRgbLong $officeColors[$i - 1]
}
